$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H1's value: PLATFORM description string with expanded platform list
$ws.Range("H1").Value = "PLATFORM:String:Default`n(CISCO_IOS|CISCO_XE|CISCO_NXOS|CISCO_WLC_CAT|CISCO_XR)"

# Touch the alignment property on H1 so a new cell style (applyAlignment) is
# recorded for it, matching the second cellXfs entry introduced upstream.
$ws.Range("H1").IndentLevel = 0

# Move the active selection to H1
$ws.Range("H1").Select()
